$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.726.19"

$ws.Range("D3").Value = "2.075.88"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'233.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").Value = "'0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'58.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").Value = "'0.392"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "'0.0784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").Value = "2.382.26"
$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("D13").Value = "'14.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("D14").Value = "'20.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").Value = "'0.773"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").Value = "'5.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("D17").Value = "2.062.34"
$ws.Range("E17").Value = "  -2.00%  "

$ws.Range("D18").Value = "37.684.28"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").Value = "'6.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").Value = "'71.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  +1.35%  "

$ws.Range("D22").Value = "'228.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("E24").Value = "  -1.13%  "

$ws.Range("D25").Value = "'2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("D26").Value = "'169.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("E27").Value = "  +4.50%  "

$ws.Range("D28").Value = "'9.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").Value = "'19.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("D31").Value = "'0.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("D32").Value = "'4.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").Value = "'0.0630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").Value = "'4.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "

$ws.Range("D35").Value = "'2.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.43%  "

$ws.Range("D36").Value = "'1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("E37").Value = "  -3.39%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").Value = "'5.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.56%  "

$ws.Range("D40").Value = "'0.0970"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").Value = "'98.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("E42").Value = "  -2.07%  "

$ws.Range("D43").Value = "'0.0215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.57%  "

$ws.Range("D44").Value = "1.451.02"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.91%  "

$ws.Range("E46").Value = "  -0.74%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'16.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.73%  "

$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("D49").Value = "'7.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("D51").Value = "2.266.18"
$ws.Range("E51").Value = "  -1.56%  "
